# The "LIST OF ABBREVIATIONS" section had a "LIDAR  LIght Dectection and
# Ranging" entry removed entirely (whole paragraph, including its mark).
# The hidden "_GoBack" bookmark - which used to sit right before the "RF"
# entry - now sits right before the "LS" entry that slides up to take the
# deleted paragraph's place.

$d = $word.ActiveDocument

# Locate the paragraph that holds the LIDAR abbreviation entry.
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("LIDAR*Ranging", $false, $false, $true, `
    $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Host "LIDAR entry not found; nothing to do."
} else {

    # Resolve the matching Document.Paragraphs index (chaining
    # Range.Paragraphs.Item off a Find range is unreliable here, so walk
    # the document's paragraph collection and match on start offset).
    $targetIndex = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Start -eq $findRange.Start) {
            $targetIndex = $i
            break
        }
    }

    $lidarPara = $d.Paragraphs.Item($targetIndex)
    Write-Host "Removing paragraph: [$($lidarPara.Range.Text)]"

    # Delete the whole paragraph, including its trailing paragraph mark,
    # so the following paragraph (LS ...) slides up into its place.
    $lidarPara.Range.Delete()

    # The paragraph that now occupies that same index is the "LS" entry.
    $lsPara = $d.Paragraphs.Item($targetIndex)
    Write-Host "LS entry now reads: [$($lsPara.Range.Text)]"

    # Move (re-add) the hidden _GoBack bookmark to sit right before it,
    # collapsed (zero-length) at the paragraph's start, same as before.
    $d.Bookmarks.Add("_GoBack", $d.Range($lsPara.Range.Start, $lsPara.Range.Start))
}
